# Generate Report for handback
#
# The localization file "f27a5b14-744b-485f-86f2-00912247ca26.md" has been
# handed back for both the zh-cn and de-de locales, so its status moves from
# "Ready for handoff" to "Handed back: in sync with en-us" on every sheet
# that tracks it, and the per-locale "Latest Handback DateTime" is refreshed
# with the new handback timestamp.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 corresponds to f27a5b14-744b-485f-86f2-00912247ca26.md
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-us"
$overview.Range("C3").Value = "Handed back: in sync with en-us"

# zh-cn sheet: update status + latest handback datetime for the same file
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-us"
$zhcn.Range("G3").Value = "2016-01-11 03:27:37"

# de-de sheet: update status + latest handback datetime for the same file
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-us"
$dede.Range("G3").Value = "2016-01-11 03:27:58"
